$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.038291333333333
$ws.Range("H2").Value = 3.114874
$ws.Range("I2").Value = 0.09907705749147798
$ws.Range("J2").Value = 0.09907705749147799
$ws.Range("M2").Value = 200.005264
$ws.Range("N2").Value = 600.0157919999999
$ws.Range("O2").Value = 0.9465949791503665
$ws.Range("P2").Value = 0.9465949791503667
$ws.Range("Q2").Value = 207.6637322322453
$ws.Range("R2").Value = 1868.973590090208
$ws.Range("S2").Value = 0.09378584517042526
$ws.Range("T2").Value = 0.09378584517042529
$ws.Range("G3").Value = 1.038291333333333
$ws.Range("H3").Value = 3.114874
$ws.Range("I3").Value = 0.09907705749147798
$ws.Range("J3").Value = 0.09907705749147799
$ws.Range("O3").Value = 0.006425713585924051
$ws.Range("P3").Value = 0.006425713585924052
$ws.Range("Q3").Value = 1.409671184508222
$ws.Range("R3").Value = 12.687040660574
$ws.Range("S3").Value = 0.0006366407943763683
$ws.Range("T3").Value = 0.0006366407943763685
$ws.Range("G4").Value = 1.038291333333333
$ws.Range("H4").Value = 3.114874
$ws.Range("I4").Value = 0.09907705749147798
$ws.Range("J4").Value = 0.09907705749147799
$ws.Range("M4").Value = 9.926218666666667
$ws.Range("O4").Value = 0.04697930726370939
$ws.Range("P4").Value = 0.0469793072637094
$ws.Range("Q4").Value = 10.30630681437155
$ws.Range("R4").Value = 92.756761329344
$ws.Range("S4").Value = 0.004654571526676344
$ws.Range("T4").Value = 0.004654571526676345
$ws.Range("I5").Value = 0.4810466828311408
$ws.Range("J5").Value = 0.4810466828311408
$ws.Range("M5").Value = 200.005264
$ws.Range("N5").Value = 600.0157919999999
$ws.Range("O5").Value = 0.9465949791503665
$ws.Range("P5").Value = 0.9465949791503667
$ws.Range("Q5").Value = 1008.265203508373
$ws.Range("R5").Value = 9074.386831575359
$ws.Range("S5").Value = 0.4553563747048966
$ws.Range("T5").Value = 0.4553563747048968
$ws.Range("I6").Value = 0.4810466828311408
$ws.Range("J6").Value = 0.4810466828311408
$ws.Range("O6").Value = 0.006425713585924051
$ws.Range("P6").Value = 0.006425713585924052
$ws.Range("S6").Value = 0.003091068205331759
$ws.Range("T6").Value = 0.00309106820533176
$ws.Range("I7").Value = 0.4810466828311408
$ws.Range("J7").Value = 0.4810466828311408
$ws.Range("M7").Value = 9.926218666666667
$ws.Range("Q7").Value = 50.03998736760889
$ws.Range("S7").Value = 0.02259923992091232
$ws.Range("T7").Value = 0.02259923992091232
$ws.Range("I8").Value = 0.4198762596773812
$ws.Range("J8").Value = 0.4198762596773812
$ws.Range("M8").Value = 200.005264
$ws.Range("N8").Value = 600.0157919999999
$ws.Range("O8").Value = 0.9465949791503665
$ws.Range("P8").Value = 0.9465949791503667
$ws.Range("Q8").Value = 880.0530957211786
$ws.Range("R8").Value = 7920.477861490606
$ws.Range("S8").Value = 0.3974527592750446
$ws.Range("T8").Value = 0.3974527592750446
$ws.Range("I9").Value = 0.4198762596773812
$ws.Range("J9").Value = 0.4198762596773812
$ws.Range("O9").Value = 0.006425713585924051
$ws.Range("P9").Value = 0.006425713585924052
$ws.Range("S9").Value = 0.002698004586215923
$ws.Range("T9").Value = 0.002698004586215923
$ws.Range("I10").Value = 0.4198762596773812
$ws.Range("J10").Value = 0.4198762596773812
$ws.Range("M10").Value = 9.926218666666667
$ws.Range("Q10").Value = 43.67684775739378
$ws.Range("S10").Value = 0.01972549581612073
$ws.Range("T10").Value = 0.01972549581612073
